$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.879.34'
$ws.Range("E2").Value = '  -4.27%  '
$ws.Range("D3").Value = '2.429.13'
$ws.Range("E3").Value = '  -7.24%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '542.44'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.84%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.99'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.51%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.586'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -6.64%  '
$ws.Range("D9").Value = '2.428.50'
$ws.Range("E9").Value = '  -7.24%  '
$ws.Range("E10").Value = '  -10.46%  '
$ws.Range("E11").Value = '  -2.02%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.37'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -7.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.347'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -9.25%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.67'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -10.00%  '
$ws.Range("D15").Value = '2.866.17'
$ws.Range("E15").Value = '  -7.27%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000163'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -10.68%  '
$ws.Range("D17").Value = '60.735.34'
$ws.Range("E17").Value = '  -4.38%  '
$ws.Range("D18").Value = '2.429.41'
$ws.Range("E18").Value = '  -7.77%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.93'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -9.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.89'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -9.54%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.12'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -8.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '315.05'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -7.62%  '
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("E24").Value = '  -1.77%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '63.34'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.77%  '
$ws.Range("D26").Value = '2.568.20'
$ws.Range("E26").Value = '  -6.80%  '
$ws.Range("D27").Value = '0.0₃0949'
$ws.Range("E27").Value = '  -14.54%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.997'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.34%  '
$ws.Range("B29").Value = 'Fetch.AI'
$ws.Range("C29").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.45'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -12.00%  '
$ws.Range("B30").Value = 'Aptos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.64'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.10%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.13'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -11.19%  '
$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '518.06'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -9.14%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.146'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -8.97%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.89'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.12%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.56'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -9.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.62'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -14.75%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.77'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -11.63%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.373'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.93%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.24'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.98%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '142.30'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -7.43%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.74'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.83%  '
$ws.Range("E43").Value = '  +0.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.29'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.19%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.24'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -11.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '139.41'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -13.92%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.57'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -8.44%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '21.09'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -12.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0527'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -9.21%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.580'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -7.45%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0927'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.93%  '
